# IO-fix: FxE is now independent from efficiency; trd removed from being duplicated
# across FIO/FOE. Reworks the template header row (row 5):
#   - "Option" and "Period" columns are dropped
#   - "Type" and "Parameter" swap order
#   - a new "Flow" column is inserted after "Year"
#   - everything shifts left by one column (old M5 "Note" goes away, content
#     that used to spill into M now ends at L)
#   - the header row is bolded
#   - the view freezes panes below the header row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the now-unused trailing header cell (was "Note" in column M).
$ws.Range("M5").ClearContents()

# Rewrite the header row (A5 Country / B5 Entity stay the same).
$ws.Range("C5").Value = "Parameter"
$ws.Range("D5").Value = "Type"
$ws.Range("E5").Value = "Year"
$ws.Range("F5").Value = "Flow"
$ws.Range("G5").Value = "Value"
$ws.Range("H5").Value = "Unit"
$ws.Range("I5").Value = "Delete"
$ws.Range("J5").Value = "Reference"
$ws.Range("K5").Value = "Link"
$ws.Range("L5").Value = "Note"

# Bold the header row.
$ws.Range("A5:L5").Font.Bold = $true

# Freeze panes just below the header row.
$ws.Range("A6").Select()
$excel.ActiveWindow.FreezePanes = $true
